$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply updated odds values per the diff (Jogos_da_Semana_FlashScore_2024-12-11.xlsx)

# Row 3
$ws.Range("G3").Value = 1.22
$ws.Range("I3").Value = 15
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.25
$ws.Range("U3").Value = 2.63
$ws.Range("V3").Value = 1.44
$ws.Range("X3").Value = 5.5
$ws.Range("Y3").Value = 10
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 29
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 67
$ws.Range("AK3").Value = 251
$ws.Range("AN3").Value = 3
$ws.Range("AP3").Value = 21
$ws.Range("AR3").Value = 41
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 3.25
$ws.Range("AU3").Value = 12
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 12
$ws.Range("AX3").Value = 51
$ws.Range("AY3").Value = 51

# Row 4
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.83

# Row 6
$ws.Range("G6").Value = 2.4
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 3.1
$ws.Range("L6").Value = 3.4
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("W6").Value = 9
$ws.Range("Z6").Value = 23
$ws.Range("AJ6").Value = 10
$ws.Range("AL6").Value = 21
$ws.Range("AX6").Value = 15

# Row 7
$ws.Range("K7").Value = 1.95
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.53
$ws.Range("AS7").Value = 301

# Row 8
$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 3.75
$ws.Range("L8").Value = 4.5
$ws.Range("AG8").Value = 451
$ws.Range("AI8").Value = 17
$ws.Range("AW8").Value = 5.5

# Row 9
$ws.Range("N9").Value = 10
$ws.Range("Q9").Value = 2.05
$ws.Range("Y9").Value = 9.5
$ws.Range("AD9").Value = 6.5
$ws.Range("AI9").Value = 13
